$d = $word.ActiveDocument

# 1. Insert the new introductory paragraph before the first paragraph.
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs(1)
$newParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Make sure you have Atlas on your computer.  Checkout the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>org.mozilla.rhino_x.x.x</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>plugin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from the Zeus SVN, then copy and paste it into your Atlas’ eclipse/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dropins</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> folder.</w:t></w:r></w:p>
'@
$newPara.Range.InsertXML($newParaXml)

# Picture 1 paragraph: add lastRenderedPageBreak
$p3 = $d.Paragraphs(3)
$xml3 = @'
<w:p w:rsidR="00505CB1" w:rsidRDefault="00E1671B"><w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="5943600" cy="3714750"/><wp:effectExtent l="19050" t="0" r="0" b="0"/><wp:docPr id="1" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId4" cstate="print"/><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5943600" cy="3714750"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:noFill/><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
'@
$p3.Range.InsertXML($xml3)

# Click the Add button paragraph: remove lastRenderedPageBreak
$p6 = $d.Paragraphs(6)
$xml6 = @'
<w:p w:rsidR="00E1671B" w:rsidRDefault="00E1671B"><w:r><w:t>Click the “Add” button, and then “Next” (Leave the setting to “Nothing”)</w:t></w:r></w:p>
'@
$p6.Range.InsertXML($xml6)

# Picture 10 paragraph: add lastRenderedPageBreak
$p9 = $d.Paragraphs(9)
$xml9 = @'
<w:p w:rsidR="00E1671B" w:rsidRDefault="00E1671B"><w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="5943600" cy="3714750"/><wp:effectExtent l="19050" t="0" r="0" b="0"/><wp:docPr id="10" name="Picture 10"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 10"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId7" cstate="print"/><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5943600" cy="3714750"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:noFill/><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
'@
$p9.Range.InsertXML($xml9)

# Click Directory paragraph: remove lastRenderedPageBreak
$p10 = $d.Paragraphs(10)
$xml10 = @'
<w:p w:rsidR="00E1671B" w:rsidRDefault="00E1671B"><w:r><w:t>Click “Directory” and then “Next”</w:t></w:r></w:p>
'@
$p10.Range.InsertXML($xml10)

# Picture 16 paragraph: add lastRenderedPageBreak
$p11 = $d.Paragraphs(11)
$xml11 = @'
<w:p w:rsidR="00E1671B" w:rsidRDefault="00E1671B"><w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="5943600" cy="3714750"/><wp:effectExtent l="19050" t="0" r="0" b="0"/><wp:docPr id="16" name="Picture 16"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 16"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId8" cstate="print"/><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5943600" cy="3714750"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:noFill/><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
'@
$p11.Range.InsertXML($xml11)

# Picture 22 paragraph: remove lastRenderedPageBreak
$p14 = $d.Paragraphs(14)
$xml14 = @'
<w:p w:rsidR="00E1671B" w:rsidRDefault="00A94231"><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="5943600" cy="3714750"/><wp:effectExtent l="19050" t="0" r="0" b="0"/><wp:docPr id="22" name="Picture 22"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 22"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId10" cstate="print"/><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5943600" cy="3714750"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:noFill/><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
'@
$p14.Range.InsertXML($xml14)

# Picture 25 paragraph: add lastRenderedPageBreak
$p15 = $d.Paragraphs(15)
$xml15 = @'
<w:p w:rsidR="00A94231" w:rsidRDefault="00AA12E7"><w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="5943600" cy="3714750"/><wp:effectExtent l="19050" t="0" r="0" b="0"/><wp:docPr id="25" name="Picture 25"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 25"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId11" cstate="print"/><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5943600" cy="3714750"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:noFill/><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
'@
$p15.Range.InsertXML($xml15)

# "Also add the..." paragraph: reword + remove lastRenderedPageBreak
$p16 = $d.Paragraphs(16)
$xml16 = @'
<w:p w:rsidR="00AA12E7" w:rsidRDefault="00AA12E7"><w:r><w:t xml:space="preserve">Follow the same steps to add the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dropins</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> folder inside of that same eclipse folder (it doesn’t find it recursively apparently so it has to be added manually).  T</w:t></w:r><w:r><w:t>hen click “Finish” on the target definition.</w:t></w:r></w:p>
'@
$p16.Range.InsertXML($xml16)

Write-Host "Edit complete."
